$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.210.63'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.370.69'
$ws.Range('E3').Value = '  +2.64%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  -0.03%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.44'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.58%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.84'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.12%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.500'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +1.03%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.22'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -1.57%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0791'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.122'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +3.02%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.56'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  -3.56%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.78'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '2.734.99'
$ws.Range('E15').Value = '  +2.66%  '
$ws.Range('D16').Value = '2.347.90'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '43.200.26'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('E20').Value = '  +4.17%  '
$ws.Range('D21').Value = '0.0₃0890'
$ws.Range('E21').Value = '  -0.28%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.26'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +0.59%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.91'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  +0.14%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  -0.14%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.79'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +1.26%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.36'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -0.14%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('E30').Value = '  -2.21%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  -0.04%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.10'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('E33').Value = '  +3.88%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.32'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -1.50%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.85'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +5.17%  '
$ws.Range('E36').Value = '  -1.60%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.31'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('E39').Value = '  +2.57%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.41'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +10.53%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.109'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.36'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -37.06%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.945.63'
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('E45').Value = '  +5.03%  '
$ws.Range('E46').Value = '  -9.41%  '
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').Value = '2.592.12'
$ws.Range('E48').Value = '  +2.26%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.07'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  -0.73%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.52'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.33'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +1.05%  '
